# Add two new rows (dates 2025-11-24 and 2025-11-25, serials 45985/45986)
# to the end of each of the six worksheets in the workbook.

$wb = $excel.ActiveWorkbook

# New rows to append per worksheet (by sheet name), as (date, amount) pairs.
$newData = @{
    "한화솔루션"  = @(@(45985, 427708), @(45986, 444476))
    "아난티"      = @(@(45985, 56027),  @(45986, 57167))
    "대아티아이"  = @(@(45985, 13416),  @(45986, 13382))
    "동원산업"    = @(@(45985, 23944),  @(45986, 24424))
    "CJ씨푸드"    = @(@(45985, 3769),   @(45986, 3741))
    "사조씨푸드"  = @(@(45985, 1854),   @(45986, 1275))
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($newData.ContainsKey($name)) {
        $rows = $newData[$name]

        # Find the current last used row in column A.
        $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

        foreach ($pair in $rows) {
            $lastRow = $lastRow + 1
            $dateSerial = $pair[0]
            $amount = $pair[1]

            $dateCell = $ws.Cells.Item($lastRow, 1)
            $dateCell.Value = $dateSerial
            $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

            $amtCell = $ws.Cells.Item($lastRow, 2)
            $amtCell.Value = $amount
        }
    }
}
